$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the rows that were removed from the "QRF 2022" category
# (Bago City x2, and Kabankalan City/La Carlota City/San Carlos City/
# Silay City x2/Sipalay City). Deleting bottom-up keeps the still-pending
# top range's row numbers valid.
$ws.Range("A214:A219").EntireRow.Delete()
$ws.Range("A207:A208").EntireRow.Delete()

# Column width tweaks (C, F, R). Excel's ColumnWidth setter stores the
# value with a constant +5/6 character offset versus the raw <col width>
# seen in the XML, so subtract it to land on the exact target widths.
$offset = 5/6
$ws.Columns.Item(3).ColumnWidth = 12 - $offset
$ws.Columns.Item(6).ColumnWidth = 18 - $offset
$ws.Columns.Item(18).ColumnWidth = 42 - $offset
